# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$changes = @{
    8  = 78
    14 = 289
    15 = 24
    16 = 351
    22 = 870
    23 = 1385
    25 = 318
    26 = 185
    27 = 66
    35 = 48
    40 = 3499
    43 = 889
    45 = 61
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $changes.Keys) {
        $ws.Range("F$row").Value = $changes[$row]
    }
}
